$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # CVE-ICS_ATT&CK-Mapping
$ws2 = $wb.Worksheets.Item(2)   # ICS_ATT&CK-Intel_TAL-Mapping
$ws3 = $wb.Worksheets.Item(3)   # Mitigations-IEC_62443-Mapping

# --- Knowledge base correction on the ICS_ATT&CK-Intel_TAL-Mapping sheet ---
# MITRE renamed technique T0830 from "Man in the Middle" to "Adversary-in-the-Middle".
$ws2.Range("B16").Value = "Adversary-in-the-Middle"

# Add the two techniques referenced elsewhere in the workbook (Rogue Master /
# Spearphishing Attachment) that were missing from the TAL mapping table.
# Copy formatting from the row above down into the two new rows first.
$ws2.Range("A19:D19").Copy()
$ws2.Range("A20:D20").PasteSpecial(-4122)
$ws2.Range("A21:D21").PasteSpecial(-4122)

$ws2.Range("A20").Value = 19
$ws2.Range("B20").Value = "Rogue Master"
$ws2.Range("C20").Value = "Adept"
$ws2.Range("D20").Value = "Organization"

$ws2.Range("A21").Value = 20
$ws2.Range("B21").Value = "Spearphishing Attachment"
$ws2.Range("C21").Value = "Adept"
$ws2.Range("D21").Value = "Organization"

# Hyperlink the two newly added ICS techniques to their MITRE ATT&CK pages
# (same convention used by the existing rows in this table).
$ws2.Hyperlinks.Add($ws2.Range("B20"), "https://attack.mitre.org/techniques/T0848/")
$ws2.Hyperlinks.Add($ws2.Range("B21"), "https://attack.mitre.org/techniques/T0865/")

# --- Paper / selection update: re-point the active tab and remembered
# selections on each sheet ---
$ws1.Activate()
$ws1.Range("D34").Select()

$ws3.Activate()
$ws3.Range("F16").Select()

$ws2.Activate()
$ws2.Range("C21:D21").Select()
